# Updated cryptos list on Thu Dec 28 18:27:09 UTC 2023 with GitHub Actions
#
# Refreshes the live Price (column D) and Volume(1h) (column E) figures for
# each coin row, and swaps out one roster entry (Cronos -> Maker) in row 51.
#
# NOTE: several of the new "Price" strings (e.g. "329.76", "0.636") look
# like plain decimal numbers to Excel's automatic type detection, whereas
# the workbook stores them as plain text (matching values like "42.648.96"
# that contain more than one '.' and can only be text). To keep those
# cells text-typed - exactly like the rest of the column - we mark them
# with a text number format immediately before writing the value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.648.96'
$ws.Range("E2").Value = '  -1.10%  '

$ws.Range("D3").Value = '2.359.28'
$ws.Range("E3").Value = '  +0.52%  '

$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '329.76'
$ws.Range("E5").Value = '  +5.88%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.14'
$ws.Range("E6").Value = '  -8.12%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.636'
$ws.Range("E7").Value = '  -0.50%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.632'
$ws.Range("E9").Value = '  +0.30%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.62'
$ws.Range("E10").Value = '  -7.36%  '

$ws.Range("E11").Value = '  -1.81%  '

$ws.Range("E12").Value = '  -5.44%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.01'
$ws.Range("E13").Value = '  -4.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.106'
$ws.Range("E14").Value = '  +0.17%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.33'
$ws.Range("E15").Value = '  +0.70%  '

$ws.Range("D16").Value = '2.707.67'
$ws.Range("E16").Value = '  +0.39%  '

$ws.Range("D17").Value = '2.354.52'
$ws.Range("E17").Value = '  +0.16%  '

$ws.Range("D18").Value = '42.546.28'
$ws.Range("E18").Value = '  -1.24%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.83'
$ws.Range("E19").Value = '  +7.42%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000106'
$ws.Range("E20").Value = '  -2.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '3.74'
$ws.Range("E21").Value = '  +9.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.61'
$ws.Range("E22").Value = '  +0.62%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '269.68'
$ws.Range("E23").Value = '  +6.99%  '

$ws.Range("E24").Value = '  -9.92%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.90'
$ws.Range("E25").Value = '  +9.55%  '

$ws.Range("E26").Value = '  -0.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.45'
$ws.Range("E27").Value = '  -4.83%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.29'
$ws.Range("E28").Value = '  +3.34%  '

$ws.Range("E29").Value = '  -3.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '175.84'
$ws.Range("E30").Value = '  +1.28%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.08'
$ws.Range("E31").Value = '  -2.60%  '

$ws.Range("E32").Value = '  -2.34%  '

$ws.Range("E33").Value = '  -9.59%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.03'

$ws.Range("E35").Value = '  -0.44%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.60'
$ws.Range("E36").Value = '  -7.98%  '

$ws.Range("E37").Value = '  -4.94%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.92'
$ws.Range("E38").Value = '  +8.23%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.106'
$ws.Range("E39").Value = '  +1.58%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.82'
$ws.Range("E40").Value = '  -6.82%  '

$ws.Range("E41").Value = '  +2.36%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.233'
$ws.Range("E42").Value = '  +0.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.24'
$ws.Range("E43").Value = '  -3.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '117.74'
$ws.Range("E45").Value = '  +7.48%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.57'
$ws.Range("E46").Value = '  +29.48%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.92'
$ws.Range("E47").Value = '  -6.47%  '

$ws.Range("E48").Value = '  -2.91%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.12'
$ws.Range("E49").Value = '  -2.44%  '

$ws.Range("E50").Value = '  -2.65%  '

$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '1.568.19'
$ws.Range("E51").Value = '  +5.08%  '
